$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 18-22 (sheet shrinks from A1:E22 to A1:E17)
$ws.Range("A18:E22").Delete() | Out-Null

# Update rows 2-17 with new label (col A) and statistic values (cols B-E)
$ws.Cells.Item(2, 1).Value = "CyclomaticComplexity(CC) & NbOperators"
$ws.Cells.Item(2, 2).Value = 45
$ws.Cells.Item(2, 3).Value = 512.000000
$ws.Cells.Item(2, 4).Value = 0.062082
$ws.Cells.Item(2, 5).Value = 0.950498

$ws.Cells.Item(3, 1).Value = "CyclomaticComplexity(CC) & EffortToImplement"
$ws.Cells.Item(3, 2).Value = 48
$ws.Cells.Item(3, 3).Value = 433.000000
$ws.Cells.Item(3, 4).Value = 1.589764
$ws.Cells.Item(3, 5).Value = 0.111889

$ws.Cells.Item(4, 1).Value = "MaintainabilityIndex & MaintainabilityIndex"
$ws.Cells.Item(4, 2).Value = 41
$ws.Cells.Item(4, 3).Value = 312.000000
$ws.Cells.Item(4, 4).Value = 1.535566
$ws.Cells.Item(4, 5).Value = 0.124646

$ws.Cells.Item(5, 1).Value = "MaintainabilityIndex & ProgramVolume"
$ws.Cells.Item(5, 2).Value = 48
$ws.Cells.Item(5, 3).Value = 403.000000
$ws.Cells.Item(5, 4).Value = 1.897461
$ws.Cells.Item(5, 5).Value = 0.057768

$ws.Cells.Item(6, 1).Value = "NbUniqueOperators & NbUniqueOperators"
$ws.Cells.Item(6, 2).Value = 44
$ws.Cells.Item(6, 3).Value = 343.000000
$ws.Cells.Item(6, 4).Value = 1.773869
$ws.Cells.Item(6, 5).Value = 0.076086

$ws.Cells.Item(7, 1).Value = "NbOperators & CyclomaticComplexity(CC)"
$ws.Cells.Item(7, 2).Value = 47
$ws.Cells.Item(7, 3).Value = 517.500000
$ws.Cells.Item(7, 4).Value = 0.492070
$ws.Cells.Item(7, 5).Value = 0.622670

$ws.Cells.Item(8, 1).Value = "NbOperators & EffortToImplement"
$ws.Cells.Item(8, 2).Value = 48
$ws.Cells.Item(8, 3).Value = 422.000000
$ws.Cells.Item(8, 4).Value = 1.702586
$ws.Cells.Item(8, 5).Value = 0.088646

$ws.Cells.Item(9, 1).Value = "ProgramVolume & MaintainabilityIndex"
$ws.Cells.Item(9, 2).Value = 48
$ws.Cells.Item(9, 3).Value = 413.000000
$ws.Cells.Item(9, 4).Value = 1.794895
$ws.Cells.Item(9, 5).Value = 0.072671

$ws.Cells.Item(10, 1).Value = "DifficultyLevel & DifficultyLevel"
$ws.Cells.Item(10, 2).Value = 4
$ws.Cells.Item(10, 3).Value = 3.000000
$ws.Cells.Item(10, 4).Value = 0.730297
$ws.Cells.Item(10, 5).Value = 0.465209

$ws.Cells.Item(11, 1).Value = "DifficultyLevel & TimeToImplement"
$ws.Cells.Item(11, 2).Value = 27
$ws.Cells.Item(11, 3).Value = 114.000000
$ws.Cells.Item(11, 4).Value = 1.801875
$ws.Cells.Item(11, 5).Value = 0.071566

$ws.Cells.Item(12, 1).Value = "ProgramLevel & ProgramLevel"
$ws.Cells.Item(12, 2).Value = 3
$ws.Cells.Item(12, 3).Value = 0.000000
$ws.Cells.Item(12, 4).Value = 1.603567
$ws.Cells.Item(12, 5).Value = 0.108810

$ws.Cells.Item(13, 1).Value = "EffortToImplement & CyclomaticComplexity(CC)"
$ws.Cells.Item(13, 2).Value = 47
$ws.Cells.Item(13, 3).Value = 431.000000
$ws.Cells.Item(13, 4).Value = 1.407427
$ws.Cells.Item(13, 5).Value = 0.159302

$ws.Cells.Item(14, 1).Value = "EffortToImplement & NbOperators"
$ws.Cells.Item(14, 2).Value = 48
$ws.Cells.Item(14, 3).Value = 424.000000
$ws.Cells.Item(14, 4).Value = 1.682073
$ws.Cells.Item(14, 5).Value = 0.092556

$ws.Cells.Item(15, 1).Value = "EffortToImplement & EffortToImplement"
$ws.Cells.Item(15, 2).Value = 4
$ws.Cells.Item(15, 3).Value = 3.000000
$ws.Cells.Item(15, 4).Value = 0.730297
$ws.Cells.Item(15, 5).Value = 0.465209

$ws.Cells.Item(16, 1).Value = "TimeToImplement & DifficultyLevel"
$ws.Cells.Item(16, 2).Value = 27
$ws.Cells.Item(16, 3).Value = 114.000000
$ws.Cells.Item(16, 4).Value = 1.801875
$ws.Cells.Item(16, 5).Value = 0.071566

$ws.Cells.Item(17, 1).Value = "TimeToImplement & TimeToImplement"
$ws.Cells.Item(17, 2).Value = 4
$ws.Cells.Item(17, 3).Value = 3.000000
$ws.Cells.Item(17, 4).Value = 0.730297
$ws.Cells.Item(17, 5).Value = 0.465209
